# Reset the "Bookings" worksheet back into a blank template:
#  - row 1 (headers) is left untouched
#  - row 2 becomes a placeholder/example row ("xxxxxxxxx" / "CounselorName" / ...)
#  - rows 3-12 are wiped (date/number columns keep their number formatting,
#    the rest of the row is fully cleared)
#  - the active selection moves to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> template/placeholder values
$ws.Range("A2").Value = "xxxxxxxxx"
$ws.Range("B2").Value = "Career Exploration and Planning"
$ws.Range("C2").Value = "Virtual"
$ws.Range("E2").Value = "10:00 AM PT"
$ws.Range("G2").Value = "CounselorName"
# D2 (date) and F2 (length) already hold the values the template keeps, so
# they are left as-is.

# Rows 3-12 -> wipe out old sample bookings.
# Columns A, B, C, E, F lose both their value AND formatting (Clear),
# while D (date) and G (counselor) only lose their value, keeping the
# number format / style that was already applied to those columns.
$ws.Range("A3:C12").Clear()
$ws.Range("E3:F12").Clear()
$ws.Range("D3:D12").ClearContents()
$ws.Range("G3:G12").ClearContents()

# Move the selection to E5, matching where the author last clicked.
$ws.Range("E5").Select()
